# Update balance/stability debugging values on GLOBAL RESULTS and LANDING GEARS sheets
$wb = $excel.ActiveWorkbook

$global = $wb.Worksheets.Item("GLOBAL RESULTS")
$global.Range("C2").Value = 59.75763069502311
$global.Range("C3").Value = 12.559185988371063
$global.Range("C4").Value = 23.951837906491903
$global.Range("C5").Value = 0.5505317220765564
$global.Range("C7").Value = 43.12093038555306
$global.Range("C8").Value = 12.17679231541782
$global.Range("C9").Value = 26.673077984574366
$global.Range("C10").Value = 0.6130792807323548
$global.Range("C12").Value = 43.12093038555306
$global.Range("C13").Value = 12.17679231541782
$global.Range("C14").Value = 26.673077984574366
$global.Range("C15").Value = 0.6130792807323548
$global.Range("C17").Value = 53.6519945733621
$global.Range("C18").Value = 12.41884826641613
$global.Range("C19").Value = 17.426005229321305
$global.Range("C20").Value = 0.40053580461201665
$global.Range("C22").Value = 51.66787324255058
$global.Range("C23").Value = 12.37324334301422
$global.Range("C24").Value = 24.642692233671223
$global.Range("C25").Value = 0.5664109720919794
$global.Range("C27").Value = 0.29186955894857564
$global.Range("C28").Value = 0.6562640112814851

$landingGears = $wb.Worksheets.Item("LANDING GEARS")
$landingGears.Range("C2").Value = 12.318240617784834
